$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column C ("Section_Name"), shifting old C..I to D..J
$ws.Columns("C").Insert()

# Header row
$ws.Range("C1").Value = 'Section_Name'

# Row 2: Divyansh Bhatt / COD1
$ws.Range("A2").Value = 'Divyansh Bhatt'
$ws.Range("B2").Value = 'divyansh.bhatt@ltimindtree.com'
$ws.Range("C2").Value = 'COD1'
$ws.Range("D2").Value = 18.899999999999995
$ws.Range("E2").Value = 30
$ws.Range("F2").Value = '2025-09-09 | 01:17:41 PM'
$ws.Range("G2").ClearContents()
$ws.Range("H2").ClearContents()
$ws.Range("I2").Value = 'Based on the logs and description provided, the Contact Management System seems to have issues with the `AddContact` method, which is not handling unique contact IDs correctly, and the `DisplayContacts` method, which is not displaying contact details accurately. Additionally, the system is not handling scenarios where no contacts are available, resulting in incorrect output.'
$ws.Range("J2").Value = 'https://admin.ltimindtree.iamneo.ai/result?testId=U2FsdGVkX186B78%2FbFnb9OxUmKb67U5QNfQVK0Y2LdBOKmYPtJ1Kfp4dtMoscu%2F9VzM0lL1T%2BVBd2ad%2BmBYVAzanisG4B5HRdsa8Zkdxf2ajSkd5fxR3AAO73M%2B7J%2FZf%2BEgsxsrYl1FKqmypEgMqsQ%3D%3D'

# Row 3: Divyansh Bhatt / COD2
$ws.Range("A3").Value = 'Divyansh Bhatt'
$ws.Range("B3").Value = 'divyansh.bhatt@ltimindtree.com'
$ws.Range("C3").Value = 'COD2'
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 30
$ws.Range("F3").Value = '2025-09-09 | 01:17:41 PM'
$ws.Range("G3").ClearContents()
$ws.Range("H3").ClearContents()
$ws.Range("I3").Value = 'Based on the logs and description provided, it appears that there are issues with the implementation of the Vehicle Management System, specifically with the `AddVehicle`, `ListVehicles`, and `DeleteVehicle` methods. The logs suggest that there are discrepancies in the expected and actual outputs, indicating errors in the handling of vehicle data, such as duplicate vehicle IDs, incorrect data display, and improper error messages. 
Here is the analysis:
**AddVehicle Method:** 
The `AddVehicle` method is not working as expected, as it is throwing an error message "A contact with ID1 already exists" instead of "Vehicle added successfully.". 
**ListVehicles Method:** 
The `ListVehicles` method is not displaying the vehicles correctly, as there is a mismatch in the expected and actual outputs.
**DeleteVehicle Method:** 
There is no log provided for `DeleteVehicle` method but as per description it need to be handled if vehicle id not found then "Vehicle not found" message. 
Overall, these issues indicate that the Vehicle Management System requires further debugging and testing to ensure that it functions correctly and meets the requirements specified in the description.'
$ws.Range("J3").Value = 'https://admin.ltimindtree.iamneo.ai/result?testId=U2FsdGVkX186B78%2FbFnb9OxUmKb67U5QNfQVK0Y2LdBOKmYPtJ1Kfp4dtMoscu%2F9VzM0lL1T%2BVBd2ad%2BmBYVAzanisG4B5HRdsa8Zkdxf2ajSkd5fxR3AAO73M%2B7J%2FZf%2BEgsxsrYl1FKqmypEgMqsQ%3D%3D'

# Row 4 (new row): Rohith Kumar Thodeti / COD1
$ws.Range("A4").Value = 'Rohith Kumar Thodeti'
$ws.Range("B4").Value = 'thodeti.rohithkumar@ltimindtree.com'
$ws.Range("C4").Value = 'COD1'
$ws.Range("D4").Value = 16.5
$ws.Range("E4").Value = 30
$ws.Range("F4").Value = '2025-06-23 | 03:58:42 PM'
$ws.Range("I4").Value = 'Based on the logs and description provided, the Apartment Management System implementation seems to have several issues with its methods, including `DisplayApartmentDetails`, `MarkAsRented`, `SearchApartment`, `UpdateApartment`, and `DisplayApartments`. The logs indicate that the implementation does not match the expected output and behavior, suggesting problems with the logic and formatting of the methods. Overall, the system requires corrections to ensure it functions as intended.'
$ws.Range("J4").Value = 'https://admin.ltimindtree.iamneo.ai/result?testId=U2FsdGVkX19LZo%2F0oePpL9zuP4tV0CGFTjBquo0oBxckMTFU2G8UmHU86MsPy3GZi7v4YVIjbqv2Kq%2BaIQtRKRbHFNCylway9hgDEy1ntYUuXEhFRVYTlrzPtMmPMIDkia0gVQjhiq%2B6xLGLbYtecQ%3D%3D'
